$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 27 de Junio de 2020 a las 23:56'
$ws.Range("B4").Value = 2590266
$ws.Range("C4").Value = 37310
$ws.Range("D4").Value = 1078326
$ws.Range("E4").Value = 1383832
$ws.Range("G4").Value = 468
$ws.Range("H4").Value = 128108
$ws.Range("B5").Value = 1313667
$ws.Range("C5").Value = 33613
$ws.Range("E5").Value = 559071
$ws.Range("G5").Value = 961
$ws.Range("H5").Value = 57070
$ws.Range("B21").Value = 131800
$ws.Range("C21").Value = 7210
$ws.Range("D21").Value = 67094
$ws.Range("E21").Value = 62293
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = 2413
$ws.Range("B31").Value = 54574
$ws.Range("C31").Value = 718
$ws.Range("D31").Value = 26920
$ws.Range("E31").Value = 23230
$ws.Range("D43").Value = 29100
$ws.Range("E43").Value = 493
$ws.Range("B50").Value = 25267
$ws.Range("C50").Value = 462
$ws.Range("D50").Value = 19781
$ws.Range("E50").Value = 5408
$ws.Range("E54").Value = 7415
$ws.Range("G54").Value = 16
$ws.Range("H54").Value = 166
$ws.Range("A113").Value = 'Paraguay'
$ws.Range("B113").Value = 1942
$ws.Range("C113").Value = 231
$ws.Range("D113").Value = 1045
$ws.Range("E113").Value = 882
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 15
$ws.Range("A114").Value = 'Islandia'
$ws.Range("B114").Value = 1836
$ws.Range("C114").Value = 4
$ws.Range("D114").Value = 1814
$ws.Range("E114").Value = 12
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 10
$ws.Range("A115").Value = 'Estado de Palestina'
$ws.Range("B115").Value = 1815
$ws.Range("C115").Value = 258
$ws.Range("D115").Value = 446
$ws.Range("E115").Value = 1365
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 4
$ws.Range("A116").Value = 'Lituania'
$ws.Range("B116").Value = 1813
$ws.Range("C116").Value = 5
$ws.Range("D116").Value = 1503
$ws.Range("E116").Value = 232
$ws.Range("H116").Value = 78
$ws.Range("A117").Value = 'Libano'
$ws.Range("B117").Value = 1719
$ws.Range("C117").Value = 22
$ws.Range("D117").Value = 1153
$ws.Range("E117").Value = 533
$ws.Range("H117").Value = 33
$ws.Range("B146").Value = 713
$ws.Range("C146").Value = 1
$ws.Range("E146").Value = 481
$ws.Range("A163").Value = 'Angola'
$ws.Range("B163").Value = 259
$ws.Range("C163").Value = 47
$ws.Range("D163").Value = 81
$ws.Range("E163").Value = 168
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 10
$ws.Range("A164").Value = 'Siria'
$ws.Range("B164").Value = 256
$ws.Range("C164").Value = 1
$ws.Range("D164").Value = 102
$ws.Range("E164").Value = 145
$ws.Range("G164").Value = 1
$ws.Range("H164").Value = 9
$ws.Range("D173").Value = 134
$ws.Range("E173").Value = 3
$ws.Range("A201").Value = 'Laos'
$ws.Range("A202").Value = 'Santa Lucia'
$ws.Range("A212").Value = 'Seychelles'
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0
$ws.Range("A213").Value = 'Montserrat'
$ws.Range("D213").Value = 10
$ws.Range("H213").Value = 1
